$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of data to append:
# (row, date serial "A", nuovi pos. "B", somma mobile 7gg. "C", somma mobile 7gg. per 100mila abitanti "D")
$newRows = @(
    @{ Row = 245; A = 44319; B = 0; C = 17; D = 106.4495929868503 },
    @{ Row = 246; A = 44320; B = 2; C = 17; D = 106.4495929868503 },
    @{ Row = 247; A = 44321; B = 0; C = 15; D = 93.9261114589856 }
)

foreach ($r in $newRows) {
    $rowIndex = $r.Row

    # Column A carries the same formatting (date style, border, alignment) as the row above it
    $aboveA = $ws.Cells.Item($rowIndex - 1, 1)
    $cellA = $ws.Cells.Item($rowIndex, 1)
    $cellA.Value = $r.A
    $aboveA.Copy()
    $cellA.PasteSpecial(-4122)

    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
    $ws.Cells.Item($rowIndex, 4).Value = $r.D
}

$excel.CutCopyMode = 0
